# Append: 2025-09-28 06:24 JST
# The scraper re-ran and produced a fresh top-N list. The old rows 2-4 are
# overwritten in place with the new top-3 entries, and the old rows 5-11
# (which fell off the new top list) are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Drop all existing hyperlinks up front -- rows 5-11 are going away and
# rows 2-4 are getting brand new target URLs, so none of the old
# relationships survive as-is.
$ws.Cells.Hyperlinks.Delete()

# The new snapshot only has 3 data rows (old dimension A1:H11 -> A1:H4),
# so remove the now-unused rows 5-11 entirely (shrinks the used range).
$ws.Range("A5:H11").EntireRow.Delete()

# Column width tweaks that came with this snapshot (title column narrower,
# price column wider).
$ws.Columns.Item(2).ColumnWidth = 46.17
$ws.Columns.Item(4).ColumnWidth = 31.17

# Row 2
$ws.Range("A2").Value = "2025-09-28 06:24:50"
$ws.Range("B2").Value = "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5217096"
$ws.Range("G2").Value = 243
$ws.Range("H2").Value = "🔥API ◆ツール"

# Row 3
$ws.Range("A3").Value = "2025-09-28 06:24:50"
$ws.Range("B3").Value = "【急募】ストレスチェックサービスの開発をお手伝いください!"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5402038"
$ws.Range("G3").Value = 75
$ws.Range("H3").Value = "◆開発"

# Row 4 (no skill-summary cell this time, same as before)
$ws.Range("A4").Value = "2025-09-28 06:24:50"
$ws.Range("B4").Value = "限定公開 PR 限定公開の仕事"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5399347"
$ws.Range("G4").Value = 13

# Re-create the hyperlinks on the URL column for the surviving rows,
# pointing at the new lancers.jp detail pages.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5217096")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5402038")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5399347")

# Hyperlinks.Add() re-points the cell style to a freshly minted xf; put the
# original "Hyperlink" cell style back so F2:F4 look like they did before.
$ws.Range("F2:F4").Style = "Hyperlink"
